$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @("TAO-USD", "IMX-USD", "GRT-USD", "PEPE-USD", "MNT-USD")

$startRow = 409
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
